$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "8:25-8:30"
$ws.Range("C3").Value = "8:30-8:35"

$ws.Range("C10").Select()
